$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (268) down to the new rows (269-272)
$ws.Range("A268:R268").Copy()
$ws.Range("A269:R272").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 269
$ws.Range("A269").Value = 267
$ws.Range("B269").Value = 44509
$ws.Range("C269").Value = 1173.599975585938
$ws.Range("D269").Value = 1174.5
$ws.Range("E269").Value = 1011.52001953125
$ws.Range("F269").Value = 1023.5
$ws.Range("G269").Value = 1023.5
$ws.Range("H269").Value = 59105800
$ws.Range("I269").Value = "TSLA"
$ws.Range("J269").Value = -139.43994140625
$ws.Range("K269").Value = 25.78321402413505
$ws.Range("L269").Value = 9.484287806919642
$ws.Range("M269").Value = 73.10757123554387
$ws.Range("N269").Value = 597195900
$ws.Range("O269").Value = -14329750
$ws.Range("P269").Value = 7427049.543542396
$ws.Range("Q269").Value = -0.989013671875
$ws.Range("R269").Value = 11.13884257077492

# Row 270
$ws.Range("A270").Value = 268
$ws.Range("B270").Value = 44510
$ws.Range("C270").Value = 1010.409973144531
$ws.Range("D270").Value = 1078.099975585938
$ws.Range("E270").Value = 987.3099975585938
$ws.Range("F270").Value = 1067.949951171875
$ws.Range("G270").Value = 1067.949951171875
$ws.Range("H270").Value = 42802700
$ws.Range("I270").Value = "TSLA"
$ws.Range("J270").Value = 44.449951171875
$ws.Range("K270").Value = 24.72142900739397
$ws.Range("L270").Value = 12.87678745814732
$ws.Range("M270").Value = 65.75160029213387
$ws.Range("N270").Value = 554393200
$ws.Range("O270").Value = 10301420
$ws.Range("P270").Value = 11143725.21198066
$ws.Range("Q270").Value = -44.76900634765625
$ws.Range("R270").Value = 17.97653308323196

# Row 271
$ws.Range("A271").Value = 269
$ws.Range("B271").Value = 44511
$ws.Range("C271").Value = 1102.77001953125
$ws.Range("D271").Value = 1104.969970703125
$ws.Range("E271").Value = 1054.680053710938
$ws.Range("F271").Value = 1063.510009765625
$ws.Range("G271").Value = 1063.510009765625
$ws.Range("H271").Value = 22396600
$ws.Range("I271").Value = "TSLA"
$ws.Range("J271").Value = -4.43994140625
$ws.Range("K271").Value = 24.58321489606585
$ws.Range("L271").Value = 12.87678745814732
$ws.Range("M271").Value = 65.62523585453258
$ws.Range("N271").Value = 531996600
$ws.Range("O271").Value = 14888930
$ws.Range("P271").Value = 9346524.592201816
$ws.Range("Q271").Value = -52.25101318359376
$ws.Range("R271").Value = 15.05813982522095

# Row 272
$ws.Range("A272").Value = 270
$ws.Range("B272").Value = 44512
$ws.Range("C272").Value = 1047.5
$ws.Range("D272").Value = 1054.5
$ws.Range("E272").Value = 1019.200012207031
$ws.Range("F272").Value = 1033.420043945312
$ws.Range("G272").Value = 1033.420043945312
$ws.Range("H272").Value = 25182300
$ws.Range("I272").Value = "TSLA"
$ws.Range("J272").Value = -30.0899658203125
$ws.Range("K272").Value = 19.9096439906529
$ws.Range("L272").Value = 14.10999843052455
$ws.Range("M272").Value = 58.52396607866465
$ws.Range("N272").Value = 557178900
$ws.Range("O272").Value = 7100750
$ws.Range("P272").Value = 11766475.92774631
$ws.Range("Q272").Value = -41.214990234375
$ws.Range("R272").Value = 18.04787973148568

$ws.Range("A1").Select()
